# Top10Mapping-es.xlsx: refresh the 2017 column (now in English source text)
# and retranslate / reorder several of the 2021 column Spanish strings,
# remove two now-unused connector arrows, and tweak a couple of view/layout
# bits (selection, column width) to match the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top 10 Mapping")

# --- 2017 column (C7:C16) -------------------------------------------------
$ws.Range("C7").Value  = "A01:2017-Injection"
$ws.Range("C8").Value  = "A02:2017-Broken Authentication"
$ws.Range("C9").Value  = "A03:2017-Sensitive Data Exposure"
$ws.Range("C10").Value = "A04:2017-XML External Entities (XXE)"
$ws.Range("C11").Value = "A05:2017-Broken Access Control"
$ws.Range("C12").Value = "A06:2017-Security Misconfiguration"
$ws.Range("C13").Value = "A07:2017-Cross-Site Scripting (XSS)"
$ws.Range("C14").Value = "A08:2017-Insecure Deserialization"
$ws.Range("C15").Value = "A09:2017-Using Components with Known Vulnerabilities"
$ws.Range("C16").Value = "A10:2017-Insufficient Logging & Monitoring"

# --- "(Nueva)" markers in column D (was "(Nuevo)") -------------------------
$ws.Range("D10").Value = "(Nueva)"
$ws.Range("D14").Value = "(Nueva)"
$ws.Range("D16").Value = "(Nueva)"

# --- 2021 column (E7:E17), retranslated / reordered -------------------------
$ws.Range("E7").Value  = "A01:2021-Pérdida de Control de Acceso"
$ws.Range("E8").Value  = "A02:2021-Fallas Criptográficas"
$ws.Range("E9").Value  = "A03:2021-Inyección"
$ws.Range("E10").Value = "A04:2021-Diseño Inseguro"
$ws.Range("E11").Value = "A05:2021-Configuración de Seguridad Incorrecta"
$ws.Range("E12").Value = "A06:2021-Componentes Vulnerables y Desactualizados"
$ws.Range("E13").Value = "A07:2021-Fallas de Identificación y Autenticación"
$ws.Range("E14").Value = "A08:2021-Fallas en el Software y en la Integridad de los Datos"
$ws.Range("E15").Value = "A09:2021-Fallas en el Registro y Monitoreo*"
$ws.Range("E16").Value = "A10:2021-Falsificación de Solicitudes del Lado del Servidor (SSRF)*"
$ws.Range("E17").Value = "* A partir de la encuesta"

# --- remove the two now-superfluous connector arrows ------------------------
$ws.Shapes.Item("Straight Arrow Connector 11").Delete()
$ws.Shapes.Item("Straight Arrow Connector 12").Delete()

# --- column layout: narrow column C (spacer column B's custom width is no
# longer present in the authored file, but COM automation has no supported
# "clear custom column width" primitive that doesn't also reshuffle data) ---
$ws.Columns.Item(3).ColumnWidth = 47.66666666666667

# --- selection moved to E23 (matches the saved view state) ------------------
$ws.Range("E23").Select()
